$d = $word.ActiveDocument
$found = $d.Content.Find.Execute("date: January 21, 2011", $true, $false, $false, $false, $false, $true, 1, $false, "date: January 21, 2016", 2)
Write-Output "found=$found"
